# Fixed update to excel issue
#
# 1. Rename "Requested quantity" header to "Weekly_PO_Qty" on "Weekly Quantity".
# 2. Rename "Requested quantity" header to "Monthly_PO_Qty" on "Monthly Trend".
# 3. Add a new "PO Forecast" sheet (after "Monthly Trend") with forecast data.

$wb = $excel.ActiveWorkbook

# --- 1 & 2: header renames -------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3: new "PO Forecast" sheet --------------------------------------------
$wsForecast = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsForecast.Name = "PO Forecast"

# Re-use the existing header/date formatting from "Weekly Quantity" so the
# new sheet matches the look of the rest of the workbook.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A39").PasteSpecial(-4122)  # xlPasteFormats

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$data = @(
  @(45039.99999999999, 109, -10.77396454417753, 227.2491891578139),
  @(45046.99999999999, 110, -16.82661761016323, 231.1432876763967),
  @(45053.99999999999, 111, -8.805329138251052, 235.6814139833062),
  @(45060.99999999999, 112, -20.40749526096801, 224.0192970578671),
  @(45067.99999999999, 113, -7.362565240712987, 238.2697645057842),
  @(45074.99999999999, 114, -8.932392449790832, 233.1294614406904),
  @(45088.99999999999, 116, 1.59572439368222, 246.4469795173899),
  @(45095.99999999999, 118, 0.1444827947117737, 240.8566901052019),
  @(45102.99999999999, 119, -5.48201626162468, 250.2114896628985),
  @(45109.99999999999, 120, -4.93889578848559, 237.6166873421571),
  @(45116.99999999999, 121, -12.8149792870835, 243.902720394745),
  @(45123.99999999999, 122, 1.411308222595669, 246.9811820815492),
  @(45130.99999999999, 123, 0.71231501715023, 248.3939189604946),
  @(45137.99999999999, 124, 11.00275314351658, 260.955582583059),
  @(45144.99999999999, 125, -2.073865380734147, 242.6397454793093),
  @(45151.99999999999, 126, 8.141600523003126, 248.374471536779),
  @(45158.99999999999, 128, 7.497678991867557, 257.9063774523412),
  @(45165.99999999999, 129, 14.73725152005726, 247.7216839771135),
  @(45172.99999999999, 130, 14.04117961126989, 262.7767217993082),
  @(45179.99999999999, 131, 4.17898173160676, 249.6758631711655),
  @(45186.99999999999, 132, 12.74257082999314, 255.4594027885977),
  @(45193.99999999999, 133, 8.31010185498169, 249.1572316701314),
  @(45200.99999999999, 134, 24.61292434401801, 263.5466394310242),
  @(45221.99999999999, 138, 13.93651487144441, 253.1215462057614),
  @(45228.99999999999, 139, 23.15898929144189, 267.5230833913658),
  @(45235.99999999999, 140, 16.26203796348327, 273.5929141654444),
  @(45242.99999999999, 141, 21.65387774440546, 269.9402996274665),
  @(45249.99999999999, 142, 8.421966295316706, 263.2187995166277),
  @(45256.99999999999, 143, 15.34618805161064, 256.9496014068784),
  @(45263.99999999999, 144, 14.2137067370449, 263.3872514432133),
  @(45270.99999999999, 146, 33.01583527533042, 273.0530499782628),
  @(45277.99999999999, 147, 30.05490012908335, 264.5096708240524),
  @(45284.99999999999, 148, 19.79306427619136, 261.2877412391096),
  @(45291.99999999999, 149, 29.74679090873617, 272.6015932038127),
  @(45298.99999999999, 150, 29.60442444370389, 270.6509648492189),
  @(45305.99999999999, 151, 24.44355469854788, 270.6448437129778),
  @(45312.99999999999, 152, 39.87699716054323, 274.3677155368827),
  @(45319.99999999999, 153, 29.77154341415017, 281.6114069133033)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $wsForecast.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
